# SF User list changes - 16 May - Initial
# Replace "Ashley Choi" with "Aadarsh Patel" in the Users sheet (row 3, column A),
# matching the border formatting already used by the rest of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# Copy the formatting already applied to B3 (border style) onto A3, then set the new name.
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Aadarsh Patel"

# Leave the selection where the author left it when saving.
$ws.Range("A8").Select()
